# Applies crypto price/volume updates per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a value to be written as text, even if it looks numeric,
# to match the original inline/shared string cell type (e.g. "1.00", "0.571").
function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "59.075.86"
$ws.Range("E2").Value = "  -3.86%  "
Set-TextValue $ws.Range("D3") "2.508.86"
$ws.Range("E3").Value = "  -1.89%  "
Set-TextValue $ws.Range("D4") "0.999"
$ws.Range("E4").Value = "  -0.13%  "
Set-TextValue $ws.Range("D5") "537.89"
$ws.Range("E5").Value = "  -1.79%  "
Set-TextValue $ws.Range("D6") "143.94"
$ws.Range("E6").Value = "  -5.08%  "
Set-TextValue $ws.Range("D8") "0.571"
$ws.Range("E8").Value = "  -2.49%  "
Set-TextValue $ws.Range("D9") "2.538.23"
$ws.Range("E9").Value = "  -0.86%  "
Set-TextValue $ws.Range("D10") "0.0999"
$ws.Range("E10").Value = "  -2.38%  "
$ws.Range("E11").Value = "  -1.81%  "
Set-TextValue $ws.Range("D12") "5.53"
$ws.Range("E12").Value = "  +1.38%  "
Set-TextValue $ws.Range("D13") "0.352"
$ws.Range("E13").Value = "  -2.80%  "
Set-TextValue $ws.Range("D14") "2.944.29"
$ws.Range("E14").Value = "  -2.19%  "
Set-TextValue $ws.Range("D15") "23.71"
$ws.Range("E15").Value = "  -5.59%  "
Set-TextValue $ws.Range("D16") "58.967.06"
$ws.Range("E16").Value = "  -3.91%  "
$ws.Range("E17").Value = "  -2.08%  "
Set-TextValue $ws.Range("D18") "2.526.06"
$ws.Range("E18").Value = "  -1.71%  "
Set-TextValue $ws.Range("D19") "11.32"
$ws.Range("E19").Value = "  -1.45%  "
$ws.Range("E20").Value = "  -4.38%  "
Set-TextValue $ws.Range("D21") "323.07"
$ws.Range("E21").Value = "  -3.83%  "
Set-TextValue $ws.Range("D22") "1.00"
$ws.Range("E22").Value = "  +0.15%  "
Set-TextValue $ws.Range("D23") "5.76"
$ws.Range("E23").Value = "  -3.20%  "
Set-TextValue $ws.Range("D24") "62.04"
$ws.Range("E24").Value = "  -1.03%  "
Set-TextValue $ws.Range("D25") "0.440"
$ws.Range("E25").Value = "  -9.05%  "
Set-TextValue $ws.Range("D26") "0.163"
$ws.Range("E26").Value = "  -1.88%  "
Set-TextValue $ws.Range("D27") "2.617.69"
$ws.Range("E27").Value = "  -2.84%  "
Set-TextValue $ws.Range("D28") "0.992"
$ws.Range("E28").Value = "  -0.83%  "
Set-TextValue $ws.Range("D29") "7.76"
$ws.Range("E29").Value = "  -3.48%  "
$ws.Range("B30").Value = "Aptos"
$ws.Range("C30").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D30") "6.74"
$ws.Range("E30").Value = "  -5.00%  "
$ws.Range("B31").Value = "PEPE"
$ws.Range("C31").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue $ws.Range("D31") "0.0₃0776"
$ws.Range("E31").Value = "  -4.85%  "
$ws.Range("E32").Value = "  -4.49%  "
Set-TextValue $ws.Range("D33") "1.21"
$ws.Range("E33").Value = "  -8.69%  "
Set-TextValue $ws.Range("D34") "0.995"
$ws.Range("E34").Value = "  -0.33%  "
Set-TextValue $ws.Range("D35") "158.92"
$ws.Range("E35").Value = "  -1.18%  "
Set-TextValue $ws.Range("D36") "1.44"
$ws.Range("E36").Value = "  +4.05%  "
Set-TextValue $ws.Range("D37") "18.57"
$ws.Range("E37").Value = "  -1.87%  "
$ws.Range("E38").Value = "  -9.07%  "
Set-TextValue $ws.Range("D39") "1.61"
$ws.Range("E39").Value = "  -8.44%  "
Set-TextValue $ws.Range("D40") "5.71"
$ws.Range("E40").Value = "  -4.13%  "
Set-TextValue $ws.Range("D41") "304.29"
$ws.Range("E41").Value = "  -5.85%  "
Set-TextValue $ws.Range("D42") "36.78"
$ws.Range("E42").Value = "  -1.04%  "
$ws.Range("E43").Value = "  -6.78%  "
Set-TextValue $ws.Range("D44") "3.67"
$ws.Range("E44").Value = "  -5.30%  "
Set-TextValue $ws.Range("D45") "0.991"
$ws.Range("E45").Value = "  -0.67%  "
Set-TextValue $ws.Range("D46") "0.601"
$ws.Range("E46").Value = "  +0.09%  "
Set-TextValue $ws.Range("D47") "10.77"
$ws.Range("E47").Value = "  -1.37%  "
Set-TextValue $ws.Range("D48") "125.67"
$ws.Range("E48").Value = "  +3.44%  "
Set-TextValue $ws.Range("D49") "0.0932"
$ws.Range("E49").Value = "  -2.71%  "
Set-TextValue $ws.Range("D50") "18.65"
$ws.Range("E50").Value = "  -2.84%  "
Set-TextValue $ws.Range("D51") "0.0518"
$ws.Range("E51").Value = "  -3.60%  "
